# Milestone 2 edit: rename Sheet1 -> Events, add Sources sheet,
# populate asset/event tables, add Notes column, add hyperlinks.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Events"

# Add the new "Sources" sheet right after "Events"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sources"

# ---- Sheet1 (Events) cell values ----
$ws1.Range("A1").Value = "Event"
$ws1.Range("B1").Value = "Description"
$ws1.Range("C1").Value = "Category"
$ws1.Range("D1").Value = "Assets Required"
$ws1.Range("E1").Value = "Status"
$ws1.Range("F1").Value = "Notes"
$ws1.Range("A2").Value = "SecurityRoomAmbience"
$ws1.Range("B2").Value = "Ambient sounds for the security room"
$ws1.Range("C2").Value = "Ambience"
$ws1.Range("D2").Value = "Several electronic beeping sounds, laser buzzing sound"
$ws1.Range("E2").Value = "Implemented"
$ws1.Range("A3").Value = "MarketAmbience"
$ws1.Range("B3").Value = "Ambient sounds for the market scene"
$ws1.Range("C3").Value = "Ambience"
$ws1.Range("D3").Value = "2-3 minimal crowd sounds"
$ws1.Range("E3").Value = "Assets Acquired"
$ws1.Range("A4").Value = "Argument"
$ws1.Range("B4").Value = "2 people arguing"
$ws1.Range("C4").Value = "Dialog"
$ws1.Range("D4").Value = "A recording of an argument"
$ws1.Range("E4").Value = "Implemented"
$ws1.Range("A5").Value = "Snoring"
$ws1.Range("B5").Value = "A snoring sound"
$ws1.Range("C5").Value = "Dialog"
$ws1.Range("D5").Value = "A snoring sound effect"
$ws1.Range("E5").Value = "Assets Acquired"
$ws1.Range("A6").Value = "DialogueAppear"
$ws1.Range("B6").Value = "Sounds for dialogue text appearing"
$ws1.Range("C6").Value = "Interface Sound"
$ws1.Range("D6").Value = "Positive, Negative sounds"
$ws1.Range("E6").Value = "Implemented"
$ws1.Range("A7").Value = "MovementInterface"
$ws1.Range("B7").Value = "Interface sound for clicking to move somewhere"
$ws1.Range("C7").Value = "Interface Sound"
$ws1.Range("D7").Value = "Simple click sound"
$ws1.Range("E7").Value = "Assets Acquired"
$ws1.Range("A8").Value = "SceneTransition"
$ws1.Range("B8").Value = "Interface sound for scene transition"
$ws1.Range("C8").Value = "Interface Sound"
$ws1.Range("D8").Value = "Woosh sound "
$ws1.Range("E8").Value = "Assets Acquired"
$ws1.Range("A9").Value = "SecurityRoomBG"
$ws1.Range("B9").Value = "Background music for the security room"
$ws1.Range("C9").Value = "Music"
$ws1.Range("D9").Value = "Security Room Background Track"
$ws1.Range("E9").Value = "Incomplete"
$ws1.Range("A10").Value = "MarketBG"
$ws1.Range("B10").Value = "Background music for the market scene"
$ws1.Range("C10").Value = "Music"
$ws1.Range("D10").Value = "Market Scene Background Track"
$ws1.Range("E10").Value = "Incomplete"
$ws1.Range("A11").Value = "Footsteps"
$ws1.Range("B11").Value = "Various footstep sounds based on ground type"
$ws1.Range("C11").Value = "Sound Effect"
$ws1.Range("D11").Value = "2-3 footstep sounds each for carpet, metal, dirt"
$ws1.Range("E11").Value = "Implemented"
$ws1.Range("A12").Value = "CollectCoin"
$ws1.Range("B12").Value = "Effect for collecting a coin"
$ws1.Range("C12").Value = "Sound Effect"
$ws1.Range("D12").Value = "Cash register sound effect"
$ws1.Range("E12").Value = "Assets Acquired"
$ws1.Range("A13").Value = "Bird"
$ws1.Range("B13").Value = "Bird chirping sound"
$ws1.Range("C13").Value = "Sound Effect"
$ws1.Range("D13").Value = "Several bird calls"
$ws1.Range("E13").Value = "Assets Acquired"
$ws1.Range("A14").Value = "GlassesDrop"
$ws1.Range("B14").Value = "Glasses hitting ground effect"
$ws1.Range("C14").Value = "Sound Effect"
$ws1.Range("D14").Value = "Sound Effect for glasses hitting ground"
$ws1.Range("E14").Value = "Assets Acquired"
$ws1.Range("A15").Value = "CollectGlasses"
$ws1.Range("B15").Value = "Effect for collecting a pair of glasses"
$ws1.Range("C15").Value = "Sound Effect"
$ws1.Range("D15").Value = "Glasses, or generic sounding picking up sound effect"
$ws1.Range("E15").Value = "Incomplete"
$ws1.Range("A16").Value = "CollectFish"
$ws1.Range("B16").Value = "Effect for collecting a fish"
$ws1.Range("C16").Value = "Sound Effect"
$ws1.Range("D16").Value = "Effect representative of a fish"
$ws1.Range("E16").Value = "Assets Acquired"
$ws1.Range("A17").Value = "CollectCoffee"
$ws1.Range("B17").Value = "Effect for collecting coffee"
$ws1.Range("C17").Value = "Sound Effect"
$ws1.Range("D17").Value = "Coffee grinding sound effect"
$ws1.Range("E17").Value = "Incomplete"
$ws1.Range("A18").Value = "Win"
$ws1.Range("B18").Value = "Effect for winning the game"
$ws1.Range("C18").Value = "Sound Effect"
$ws1.Range("D18").Value = "Win sfx"
$ws1.Range("E18").Value = "Assets Acquired"

# ---- Sheet2 (Sources) cell values ----
$ws2.Range("A1").Value = "Asset"
$ws2.Range("B1").Value = "Source"
$ws2.Range("A2").Value = "beep-*"
$ws2.Range("B2").Value = "Created from FL Studio"
$ws2.Range("A3").Value = "laser_sustained.wav"
$ws2.Range("B3").Value = "https://freesound.org/people/ledhed2222/sounds/397280/"
$ws2.Range("A4").Value = "Crowd Exterior Large Size, City Voices, Footsteps, Distant Traffic Stereo .wav"
$ws2.Range("B4").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A5").Value = "argument ambience_130611.wav"
$ws2.Range("B5").Value = "https://freesound.org/people/miastodzwiekow/sounds/122328/"
$ws2.Range("A6").Value = "MaleSnore_Raw_bip.wav"
$ws2.Range("B6").Value = "https://freesound.org/people/passAirmangrace/sounds/340893/"
$ws2.Range("A7").Value = "guitar_ac_fx_006.wav"
$ws2.Range("B7").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A8").Value = "button_002.wav"
$ws2.Range("B8").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A9").Value = "marimba_tone_007.wav"
$ws2.Range("B9").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A10").Value = "org_short_R_to_L_002.wav"
$ws2.Range("B10").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A11").Value = "S23_SFX_Footsteps_Gravel_Loafers_Loops_Walk_Normal.wav"
$ws2.Range("B11").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A12").Value = "S23_SFX_Footsteps_Metal_Boots_Loop_Jogging.wav"
$ws2.Range("B12").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A13").Value = "S23_SFX_Footsteps_Snow_Singles_Sequence.wav"
$ws2.Range("B13").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A14").Value = "Coins_Bottlecaps_Drop.wav"
$ws2.Range("B14").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A15").Value = "PRB216 Saker Falcon (Tosia) - Call, Squawk Screech Scream Cry, High Nervous - 8060 MF.wav"
$ws2.Range("B15").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A16").Value = "Glass,Plate Glass,Thick,Break,Topple,Schoeps.wav"
$ws2.Range("B16").Value = "https://sonniss.com/gameaudiogdc19/"
$ws2.Range("A17").Value = "Water_Pouring_02.wav"
$ws2.Range("B17").Value = "https://sonniss.com/gameaudiogdc19/"

# ---- Sheet2 hyperlinks (rows 3-17) ----
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://freesound.org/people/ledhed2222/sounds/397280/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B5"), "https://freesound.org/people/miastodzwiekow/sounds/122328/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B6"), "https://freesound.org/people/passAirmangrace/sounds/340893/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B7"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B8"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B9"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B10"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B11"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B12"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B13"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B14"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B15"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B16"), "https://sonniss.com/gameaudiogdc19/") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B17"), "https://sonniss.com/gameaudiogdc19/") | Out-Null

# ---- Header formatting (bold) ----
$ws1.Range("A1:F1").Font.Bold = $true
$ws2.Range("A1:B1").Font.Bold = $true

# ---- Column widths ----
$ws1.Columns.Item(6).ColumnWidth = 19.1
$ws2.Columns.Item(1).ColumnWidth = 79.1
$ws2.Columns.Item(2).ColumnWidth = 55.5

# ---- Selections / active sheet (matches saved workbook state) ----
$ws2.Range("B19").Select()
$ws1.Activate()
$ws1.Range("E7").Select()
